$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Start from a clean sheet and re-lay out all content so that stale cells
# left behind by the row reflow (old rows 5-9, 13) are removed.
$ws.Cells.Clear()

# Helper-free, explicit cell-by-cell population matching the target layout.

# --- Name / title block -------------------------------------------------
$ws.Range("A1").Font.Size = 18
$ws.Range("A1").Value = "Uzbekistan"

$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Value = "MSME Participation on the Economy"

$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Value = "Source Type: Statistical Institution (Most Widely Used)"

# --- First table: MSME participation ------------------------------------
$ws.Range("B9").Font.Bold = $true
$ws.Range("B9").Value = "Micro"
$ws.Range("C9").Font.Bold = $true
$ws.Range("C9").Value = "SMEs"
$ws.Range("D9").Font.Bold = $true
$ws.Range("D9").Value = "MSMEs"

$ws.Range("A10").Font.Bold = $true
$ws.Range("A10").Value = "Enterprises (absolute #)"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "253080"

$ws.Range("A11").Font.Bold = $true
$ws.Range("A11").Value = "Enterprises density (per 1000 people)"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.5"

$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").Value = "Employment (absolute #)"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9519600"

$ws.Range("A13").Font.Italic = $true
$ws.Range("A13").Value = "Source: SCRUS, 2013"

# --- Second table: Value added -------------------------------------------
$ws.Range("B15").Font.Bold = $true
$ws.Range("B15").Value = "Micro"
$ws.Range("C15").Font.Bold = $true
$ws.Range("C15").Value = "SMEs"
$ws.Range("D15").Font.Bold = $true
$ws.Range("D15").Value = "MSMEs"

$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").Value = "Value added to the economy (% of total)"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "51"

$ws.Range("A17").Font.Italic = $true
$ws.Range("A17").Value = "Source: SCRUS, 2013"

# --- Sector distribution details -----------------------------------------
$ws.Range("A22").Font.Bold = $true
$ws.Range("A22").Value = "Sector Distribution Details"

$ws.Range("A25").Font.Bold = $true
$ws.Range("A25").Value = "SCRUS"

$longText = "The State Committee of the Republic of Uzbekistan on Statistics (SCRUS), `"По состоянию на 1 октября 2012 года количество зарегистрированных субъектов малого бизнеса (без фермерских хозяйств) составило 249 915`", 2012. Available at http://www.stat.uz/press/1/5154/?sphrase_id=108660`nThe State Committee of the Republic of Uzbekistan on Statistics (SCRUS), `"Об итогах социально-экономического развития Республики Узбекистан за I квартал 2014 года`", 2014. Available at http://www.stat.uz/press/1/8359/?sphrase_id=108660"

$ws.Range("A26").Font.Italic = $true
$ws.Range("A26").Value = $longText
